$wb = $excel.ActiveWorkbook

# Target OOXML column width of 29.9777047293527 (approx 30 widened report columns)
# and 40 are reproduced through Excel's character-width grid (1/6 quantization);
# these ColumnWidth inputs land on the nearest achievable stored width.
$wideCol = 29 + (1/6)
$fullCol = 39 + (1/6)

# --- Overview sheet: handoff status text + column widths ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = $wideCol
$wsOverview.Columns.Item(6).ColumnWidth = $wideCol

# --- zh-cn sheet: report generated for handback ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("I2").Value = "6bbff5fe-2d23-4f51-baa8-e64361157706.md"
$wsZh.Range("I2").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aefaa76cde8fe1c1b099b8d4c5fd52a83dc7492f/e2e/6bbff5fe-2d23-4f51-baa8-e64361157706.md", "", "", "6bbff5fe-2d23-4f51-baa8-e64361157706.md")
$wsZh.Range("J2").Value = "6bbff5fe-2d23-4f51-baa8-e64361157706.49c5d8b27207b496a6b4397d0ff9e440d659de58.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-15 10:56:26"
$wsZh.Columns.Item(3).ColumnWidth = $wideCol
$wsZh.Columns.Item(9).ColumnWidth = $fullCol
$wsZh.Columns.Item(10).ColumnWidth = $fullCol

# --- de-de sheet: report generated for handback ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("I2").Value = "6bbff5fe-2d23-4f51-baa8-e64361157706.md"
$wsDe.Range("I2").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/aefaa76cde8fe1c1b099b8d4c5fd52a83dc7492f/e2e/6bbff5fe-2d23-4f51-baa8-e64361157706.md", "", "", "6bbff5fe-2d23-4f51-baa8-e64361157706.md")
$wsDe.Range("J2").Value = "6bbff5fe-2d23-4f51-baa8-e64361157706.49c5d8b27207b496a6b4397d0ff9e440d659de58.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-15 10:56:33"
$wsDe.Columns.Item(3).ColumnWidth = $wideCol
$wsDe.Columns.Item(9).ColumnWidth = $fullCol
$wsDe.Columns.Item(10).ColumnWidth = $fullCol
